# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ),
# LevePrice(NQ/HQ) and LeveProfit(NQ/HQ) columns (H:N) on each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 316.8889
$ws.Range("J2").Value = 449.5
$ws.Range("L2").Value = 449.5
$ws.Range("N2").Value = -675.5
$ws.Range("H17").Value = 381.27777
$ws.Range("J17").Value = 381.27777
$ws.Range("L17").Value = 1143.83331
$ws.Range("N17").Value = -1479.83331
$ws.Range("H18").Value = 1038.5555
$ws.Range("J18").Value = 5900
$ws.Range("L18").Value = 5900
$ws.Range("N18").Value = -6468
$ws.Range("H40").Value = 6557.75
$ws.Range("I40").Value = 5074.25
$ws.Range("J40").Value = 7299.5
$ws.Range("K40").Value = 5074.25
$ws.Range("L40").Value = 7299.5
$ws.Range("M40").Value = -4899.25
$ws.Range("N40").Value = -7649.5
$ws.Range("H69").Value = 5875
$ws.Range("J69").Value = 5875
$ws.Range("L69").Value = 17625
$ws.Range("N69").Value = -19373
$ws.Range("H72").Value = 5875
$ws.Range("J72").Value = 5875
$ws.Range("L72").Value = 52875
$ws.Range("N72").Value = -61611
$ws.Range("H100").Value = 2697.6365
$ws.Range("J100").Value = 3180.5
$ws.Range("L100").Value = 3180.5
$ws.Range("N100").Value = -4262.5
$ws.Range("H137").Value = 3458.6128
$ws.Range("I137").Value = 1525.1111
$ws.Range("J137").Value = 16509.75
$ws.Range("K137").Value = 4575.3333
$ws.Range("L137").Value = 49529.25
$ws.Range("M137").Value = -2025.3333
$ws.Range("N137").Value = -54629.25
$ws.Range("H138").Value = 4796.271
$ws.Range("J138").Value = 7089.1143
$ws.Range("L138").Value = 21267.3429
$ws.Range("N138").Value = -31547.3429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1390.3529
$ws.Range("I45").Value = 874.1429000000001
$ws.Range("K45").Value = 874.1429000000001
$ws.Range("M45").Value = -497.1429000000001
$ws.Range("H61").Value = 27779938
$ws.Range("I61").Value = 33335468
$ws.Range("J61").Value = 2278.3333
$ws.Range("K61").Value = 33335468
$ws.Range("L61").Value = 2278.3333
$ws.Range("M61").Value = -33335256
$ws.Range("N61").Value = -2702.3333
$ws.Range("H74").Value = 25001454
$ws.Range("I74").Value = 33334318
$ws.Range("K74").Value = 33334318
$ws.Range("M74").Value = -33333444
$ws.Range("H77").Value = 25001454
$ws.Range("I77").Value = 33334318
$ws.Range("K77").Value = 166671590
$ws.Range("M77").Value = -166667222
$ws.Range("H110").Value = 10399.4
$ws.Range("I110").Value = 10249.25
$ws.Range("K110").Value = 10249.25
$ws.Range("M110").Value = -8204.25
$ws.Range("H132").Value = 40068340
$ws.Range("I132").Value = 14568.096
$ws.Range("J132").Value = 250350660
$ws.Range("K132").Value = 43704.288
$ws.Range("L132").Value = 751051980
$ws.Range("M132").Value = -41174.288
$ws.Range("N132").Value = -751057040
$ws.Range("H136").Value = 27779938
$ws.Range("I136").Value = 33335468
$ws.Range("J136").Value = 2278.3333
$ws.Range("K136").Value = 100006404
$ws.Range("L136").Value = 6834.999899999999
$ws.Range("M136").Value = -100003854
$ws.Range("N136").Value = -11934.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H134").Value = 3321.8948
$ws.Range("I134").Value = 3374.4707
$ws.Range("K134").Value = 10123.4121
$ws.Range("M134").Value = -7588.4121
$ws.Range("H138").Value = 67854.5
$ws.Range("I138").Value = 50709
$ws.Range("K138").Value = 50709
$ws.Range("M138").Value = -45569

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7860.4
$ws.Range("I22").Value = 9355.272000000001
$ws.Range("K22").Value = 9355.272000000001
$ws.Range("M22").Value = -9005.272000000001
$ws.Range("H31").Value = 6100779
$ws.Range("I31").Value = 2144.0303
$ws.Range("K31").Value = 2144.0303
$ws.Range("M31").Value = -1849.0303
$ws.Range("H34").Value = 6100779
$ws.Range("I34").Value = 2144.0303
$ws.Range("K34").Value = 2144.0303
$ws.Range("M34").Value = -1942.0303
$ws.Range("H41").Value = 30249.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 30249.25
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").Value = 30249.25
$ws.Range("N41").Value = -31105.25
$ws.Range("H48").Value = 37737
$ws.Range("J48").Value = 37737
$ws.Range("L48").Value = 37737
$ws.Range("N48").Value = -38689
$ws.Range("H51").Value = 19999
$ws.Range("I51").Value = 19999
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 19999
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -19263
$ws.Range("H58").Value = 1404.4546
$ws.Range("I58").Value = 1493.75
$ws.Range("J58").Value = 1166.3334
$ws.Range("K58").Value = 1493.75
$ws.Range("L58").Value = 1166.3334
$ws.Range("M58").Value = -1290.75
$ws.Range("N58").Value = -1572.3334
$ws.Range("H60").Value = 12749.75
$ws.Range("I60").Value = 3666.6667
$ws.Range("J60").Value = 39999
$ws.Range("K60").Value = 3666.6667
$ws.Range("L60").Value = 39999
$ws.Range("M60").Value = -3155.6667
$ws.Range("N60").Value = -41021
$ws.Range("H61").Value = 19999
$ws.Range("I61").Value = 19999
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 19999
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -19651
$ws.Range("H136").Value = 1404.4546
$ws.Range("I136").Value = 1493.75
$ws.Range("J136").Value = 1166.3334
$ws.Range("K136").Value = 4481.25
$ws.Range("L136").Value = 3499.0002
$ws.Range("M136").Value = -1931.25
$ws.Range("N136").Value = -8599.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1901.9584
$ws.Range("J34").Value = 4049
$ws.Range("L34").Value = 12147
$ws.Range("N34").Value = -12315
$ws.Range("H39").Value = 562
$ws.Range("I39").Value = 562
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1686
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -1392
$ws.Range("H55").Value = 696258
$ws.Range("J55").Value = 5664.143
$ws.Range("L55").Value = 16992.429
$ws.Range("N55").Value = -17346.429
$ws.Range("H98").Value = 538.7273
$ws.Range("J98").Value = 530.125
$ws.Range("L98").Value = 1590.375
$ws.Range("N98").Value = -4586.375
$ws.Range("H107").Value = 1432.6364
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1432.6364
$ws.Range("K107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("M107").Value = 4297.9092
$ws.Range("N107").Value = -8137.9092
$ws.Range("H109").Value = 4045.6924
$ws.Range("I109").Value = 6453.7144
$ws.Range("J109").Value = 1236.3334
$ws.Range("K109").Value = 19361.1432
$ws.Range("L109").Value = 3709.0002
$ws.Range("M109").Value = -18321.1432
$ws.Range("N109").Value = -5789.0002
$ws.Range("H127").Value = 2649
$ws.Range("J127").Value = 2649
$ws.Range("L127").Value = 7947
$ws.Range("N127").Value = -17867
$ws.Range("H131").Value = 1219.8966
$ws.Range("I131").Value = 737.13336
$ws.Range("J131").Value = 1737.1428
$ws.Range("K131").Value = 2211.40008
$ws.Range("L131").Value = 5211.428400000001
$ws.Range("M131").Value = 2828.59992
$ws.Range("N131").Value = -15291.4284
$ws.Range("H132").Value = 6672476.5
$ws.Range("I132").Value = 1649.75
$ws.Range("J132").Value = 11119694
$ws.Range("K132").Value = 14847.75
$ws.Range("L132").Value = 100077246
$ws.Range("M132").Value = -12317.75
$ws.Range("N132").Value = -100082306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 19538.572
$ws.Range("J49").Value = 21000
$ws.Range("L49").Value = 21000
$ws.Range("N49").Value = -21368
$ws.Range("H70").Value = 95630.73
$ws.Range("I70").Value = 137065.73
$ws.Range("J70").Value = 6841.4287
$ws.Range("K70").Value = 137065.73
$ws.Range("L70").Value = 6841.4287
$ws.Range("M70").Value = -136795.73
$ws.Range("N70").Value = -7381.4287
$ws.Range("H73").Value = 95630.73
$ws.Range("I73").Value = 137065.73
$ws.Range("J73").Value = 6841.4287
$ws.Range("K73").Value = 137065.73
$ws.Range("L73").Value = 6841.4287
$ws.Range("M73").Value = -136129.73
$ws.Range("N73").Value = -8713.4287
$ws.Range("H132").Value = 3718.4285
$ws.Range("I132").Value = 3006.75
$ws.Range("K132").Value = 9020.25
$ws.Range("M132").Value = -6490.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5420.3125
$ws.Range("I132").Value = 2673.9
$ws.Range("J132").Value = 9997.666999999999
$ws.Range("K132").Value = 8021.700000000001
$ws.Range("L132").Value = 29993.001
$ws.Range("M132").Value = -5491.700000000001
$ws.Range("N132").Value = -35053.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 20236.334
$ws.Range("J70").Value = 20236.334
$ws.Range("L70").Value = 20236.334
$ws.Range("N70").Value = -20866.334
$ws.Range("H73").Value = 20236.334
$ws.Range("J73").Value = 20236.334
$ws.Range("L73").Value = 20236.334
$ws.Range("N73").Value = -22420.334
$ws.Range("H81").Value = 8000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 0
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = 16000
$ws.Range("N81").Value = -18122
$ws.Range("H84").Value = 8000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 0
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = 80000
$ws.Range("N84").Value = -90608
$ws.Range("H100").Value = 59412948
$ws.Range("J100").Value = 1428.5714
$ws.Range("L100").Value = 2857.1428
$ws.Range("N100").Value = -3939.1428
$ws.Range("H132").Value = 1898.8055
$ws.Range("I132").Value = 1479.9062
$ws.Range("K132").Value = 4439.7186
$ws.Range("M132").Value = -1909.7186
